$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 644, shifting existing rows 644:719 down to 646:721
$ws.Rows.Item(644).Insert()
$ws.Rows.Item(644).Insert()

# New row 644 data
$ws.Range("A644").Value = 9
$ws.Range("B644").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C644").Value = "Metropolitana"
$ws.Range("D644").Value = 45212
$ws.Range("E644").Value = 13
$ws.Range("F644").Value = 100112013
$ws.Range("G644").Value = "Alcachofa"
$ws.Range("H644").Value = "Española"
$ws.Range("I644").Value = "Extra"
$ws.Range("J644").Value = 52
$ws.Range("K644").Value = 12000
$ws.Range("L644").Value = 13000
$ws.Range("M644").Value = 12500
$ws.Range("N644").Value = "$/caja 25 unidades"
$ws.Range("O644").Value = "Provincia de Limarí"
$ws.Range("P644").Value = 12500
$ws.Range("Q644").Value = 1
$ws.Range("R644").Value = "Hortaliza"

# New row 645 data
$ws.Range("A645").Value = 9
$ws.Range("B645").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C645").Value = "Metropolitana"
$ws.Range("D645").Value = 45212
$ws.Range("E645").Value = 13
$ws.Range("F645").Value = 100112013
$ws.Range("G645").Value = "Alcachofa"
$ws.Range("H645").Value = "Española"
$ws.Range("I645").Value = "Primera"
$ws.Range("J645").Value = 70
$ws.Range("K645").Value = 11000
$ws.Range("L645").Value = 12000
$ws.Range("M645").Value = 11500
$ws.Range("N645").Value = "$/caja 30 unidades"
$ws.Range("O645").Value = "Provincia de Limarí"
$ws.Range("P645").Value = 383
$ws.Range("Q645").Value = 30
$ws.Range("R645").Value = "Hortaliza"
